$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "67.089.01"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "3.103.43"
$ws.Range("E3").Value = "  -0.19%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.13%  "

Set-TextValue $ws.Range("D5") "574.01"
$ws.Range("E5").Value = "  -0.88%  "

Set-TextValue $ws.Range("D6") "177.32"
$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.102.69"
$ws.Range("E8").Value = "  -0.08%  "

Set-TextValue $ws.Range("D9") "0.512"
$ws.Range("E9").Value = "  -1.44%  "

Set-TextValue $ws.Range("D10") "6.36"
$ws.Range("E10").Value = "  -2.27%  "

Set-TextValue $ws.Range("D11") "0.152"
$ws.Range("E11").Value = "  -0.78%  "

Set-TextValue $ws.Range("D12") "0.468"
$ws.Range("E12").Value = "  -1.70%  "

$ws.Range("E13").Value = "  -2.93%  "

Set-TextValue $ws.Range("D14") "36.07"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").Value = "3.620.74"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").Value = "66.973.45"
$ws.Range("E17").Value = "  -0.27%  "

Set-TextValue $ws.Range("D18") "7.03"
$ws.Range("E18").Value = "  -0.67%  "

Set-TextValue $ws.Range("D19") "16.73"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("D20").Value = "3.100.69"
$ws.Range("E20").Value = "  -0.22%  "

Set-TextValue $ws.Range("D21") "479.17"
$ws.Range("E21").Value = "  -2.39%  "

Set-TextValue $ws.Range("D22") "7.79"
$ws.Range("E22").Value = "  -0.29%  "

Set-TextValue $ws.Range("D23") "0.688"
$ws.Range("E23").Value = "  -1.48%  "

Set-TextValue $ws.Range("D24") "83.44"
$ws.Range("E24").Value = "  -0.42%  "

Set-TextValue $ws.Range("D25") "12.59"
$ws.Range("E25").Value = "  -3.54%  "

Set-TextValue $ws.Range("D26") "2.26"
$ws.Range("E26").Value = "  -1.22%  "

Set-TextValue $ws.Range("D27") "10.13"
$ws.Range("E27").Value = "  -3.68%  "

Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  +0.00%  "

Set-TextValue $ws.Range("D29") "7.93"
$ws.Range("E29").Value = "  +0.61%  "

$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("E31").Value = "  -2.53%  "

Set-TextValue $ws.Range("D32") "28.02"
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("E33").Value = "  -1.87%  "

$ws.Range("D34").Value = "0.0₃0939"
$ws.Range("E34").Value = "  +0.19%  "

Set-TextValue $ws.Range("D35") "0.998"
$ws.Range("E35").Value = "  -0.18%  "

Set-TextValue $ws.Range("D36") "48.47"
$ws.Range("E36").Value = "  +3.72%  "

$ws.Range("E37").Value = "  -2.95%  "

Set-TextValue $ws.Range("D38") "0.943"
$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("E39").Value = "  +1.76%  "

$ws.Range("E40").Value = "  -1.88%  "

$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("E43").Value = "  -1.62%  "

Set-TextValue $ws.Range("D44") "2.72"
$ws.Range("E44").Value = "  +5.78%  "

$ws.Range("D45").Value = "2.798.44"
$ws.Range("E45").Value = "  +0.02%  "

Set-TextValue $ws.Range("D46") "372.55"
$ws.Range("E46").Value = "  -3.54%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0344"
$ws.Range("E47").Value = "  -1.60%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D48") "135.57"
$ws.Range("E48").Value = "  +0.36%  "

Set-TextValue $ws.Range("D50") "25.44"
$ws.Range("E50").Value = "  +1.94%  "

Set-TextValue $ws.Range("D51") "2.29"
$ws.Range("E51").Value = "  +4.41%  "
